$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.769.90'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.642.58'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '218.15'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('E9').Value = '  -0.63%  '
$ws.Range('D10').Value = '19.06'
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('D11').Value = '0.0847'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '1.870.80'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').Value = '1.618.50'
$ws.Range('E13').Value = '  -2.19%  '
$ws.Range('D14').Value = '4.13'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').Value = '64.64'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '26.745.21'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('E18').Value = '  -2.45%  '
$ws.Range('D20').Value = '211.11'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').Value = '4.33'
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').Value = '2.32'
$ws.Range('E23').Value = '  -5.38%  '
$ws.Range('D24').Value = '9.24'
$ws.Range('E24').Value = '  -2.50%  '
$ws.Range('D25').Value = '147.52'
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').Value = '7.06'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('E30').Value = '  -3.41%  '
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('D33').Value = '2.97'
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('D34').Value = '1.273.44'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('D38').Value = '0.526'
$ws.Range('E38').Value = '  -2.22%  '
$ws.Range('D39').Value = '0.805'
$ws.Range('E39').Value = '  -2.91%  '
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('D41').Value = '0.804'
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('D43').Value = '1.781.57'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('D45').Value = '91.38'
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').Value = '60.19'
$ws.Range('E46').Value = '  +0.87%  '
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0959'
$ws.Range('E51').Value = '  -1.58%  '
